$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 parameter values
$ws.Range("A2").Value = 100
$ws.Range("C2").Value = 10
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 200
$ws.Range("J2").Value = 100

# Remove the NumProcessors value cell in row 2 (header in L1 stays)
$ws.Range("L2").ClearContents()

# Update the active selection to A2
$ws.Range("A2").Select()
